# ==========================================================================
# edit.ps1 - apply the ruv_risks.xlsx update described by the commit
#   "update excel and upload source files"
# ==========================================================================

$wb = $excel.ActiveWorkbook

$wsData   = $wb.Worksheets.Item("data")
$wsDates  = $wb.Worksheets.Item("dates")
$wsSource = $wb.Worksheets.Item("source")

# --------------------------------------------------------------------------
# 1. "data" sheet - fill in newly-collected wave columns (U:X, plus a few
#    extra cells in already-present columns) and append three new risk rows
# --------------------------------------------------------------------------

# Row 4 - Spannungen durch Zuzug von Ausländern
$wsData.Cells.Item(4, 21).Value = 42   # U4
$wsData.Cells.Item(4, 22).Value = 41   # V4
$wsData.Cells.Item(4, 23).Value = 42   # W4
$wsData.Cells.Item(4, 24).Value = 43   # X4

# Row 7 - Vereinsamung im Alter
$wsData.Cells.Item(7, 21).Value = 39   # U7
$wsData.Cells.Item(7, 22).Value = 40   # V7
$wsData.Cells.Item(7, 23).Value = 40   # W7
$wsData.Cells.Item(7, 24).Value = 38   # X7
$wsData.Cells.Item(7, 25).Value = 37   # Y7
$wsData.Cells.Item(7, 26).Value = 46   # Z7
$wsData.Cells.Item(7, 30).Value = 32   # AD7
$wsData.Cells.Item(7, 31).Value = 35   # AE7
$wsData.Cells.Item(7, 32).Value = 38   # AF7

# Row 8 - wave/placeholder risk
$wsData.Cells.Item(8, 21).Value = 31   # U8
$wsData.Cells.Item(8, 22).Value = 29   # V8
$wsData.Cells.Item(8, 23).Value = 31   # W8
$wsData.Cells.Item(8, 24).Value = 28   # X8
$wsData.Cells.Item(8, 25).Value = 27   # Y8
$wsData.Cells.Item(8, 26).Value = 33   # Z8

# Row 10 - Drogensucht der eigenen Kinder
$wsData.Cells.Item(10, 21).Value = 38  # U10
$wsData.Cells.Item(10, 22).Value = 34  # V10
$wsData.Cells.Item(10, 23).Value = 35  # W10
$wsData.Cells.Item(10, 24).Value = 32  # X10
$wsData.Cells.Item(10, 25).Value = 31  # Y10
$wsData.Cells.Item(10, 26).Value = 36  # Z10
$wsData.Cells.Item(10, 30).Value = 24  # AD10
$wsData.Cells.Item(10, 31).Value = 19  # AE10

# Row 20 - Überforderung des Staates durch Geflüchtete
$wsData.Cells.Item(20, 26).Value = 49  # Z20

# New rows appended at the bottom of the table (28, 29, 30)
$wsData.Cells.Item(28, 1).Value = "Hinterherhinken bei Digitalisierung"
$wsData.Cells.Item(28, 31).Value = 38  # AE28

$wsData.Cells.Item(29, 1).Value = "Zerbrechen der Partnerschaft"
$wsData.Cells.Item(29, 21).Value = 18  # U29
$wsData.Cells.Item(29, 22).Value = 16  # V29
$wsData.Cells.Item(29, 23).Value = 20  # W29
$wsData.Cells.Item(29, 24).Value = 18  # X29
$wsData.Cells.Item(29, 25).Value = 15  # Y29
$wsData.Cells.Item(29, 26).Value = 21  # Z29
$wsData.Cells.Item(29, 30).Value = 10  # AD29
$wsData.Cells.Item(29, 31).Value = 15  # AE29

$wsData.Cells.Item(30, 1).Value = "häufiger Pandemien durch Globalisierung"
$wsData.Cells.Item(30, 30).Value = 42  # AD30

# widen the custom column-width run to also cover the newly used columns T:U
# (matches the on-screen width already used for columns B:S)
$wsData.Columns("T:U").ColumnWidth = $wsData.Range("S1").ColumnWidth

# --------------------------------------------------------------------------
# 2. "dates" sheet - record the exact survey start/end dates for the last
#    few waves, and refresh the respondent counts that changed
# --------------------------------------------------------------------------

# 2018 wave
$wsDates.Cells.Item(28, 3).Value = 43259   # C28
$wsDates.Cells.Item(32, 3).Copy()
$wsDates.Cells.Item(28, 3).PasteSpecial(-4122)
$wsDates.Cells.Item(28, 4).Value = 43664   # D28
$wsDates.Cells.Item(32, 4).Copy()
$wsDates.Cells.Item(28, 4).PasteSpecial(-4122)
$wsDates.Cells.Item(28, 5).Value = 2335    # E28

# 2019 wave
$wsDates.Cells.Item(29, 3).Value = 43598   # C29
$wsDates.Cells.Item(32, 3).Copy()
$wsDates.Cells.Item(29, 3).PasteSpecial(-4122)
$wsDates.Cells.Item(29, 4).Value = 43669   # D29
$wsDates.Cells.Item(32, 4).Copy()
$wsDates.Cells.Item(29, 4).PasteSpecial(-4122)
$wsDates.Cells.Item(29, 5).Value = 2446    # E29

# 2020 wave
$wsDates.Cells.Item(30, 3).Value = 43990   # C30
$wsDates.Cells.Item(32, 3).Copy()
$wsDates.Cells.Item(30, 3).PasteSpecial(-4122)
$wsDates.Cells.Item(30, 4).Value = 44033   # D30
$wsDates.Cells.Item(32, 4).Copy()
$wsDates.Cells.Item(30, 4).PasteSpecial(-4122)

# 2021 wave
$wsDates.Cells.Item(31, 3).Value = 44341   # C31
$wsDates.Cells.Item(32, 3).Copy()
$wsDates.Cells.Item(31, 3).PasteSpecial(-4122)
$wsDates.Cells.Item(31, 4).Value = 44381   # D31
$wsDates.Cells.Item(32, 4).Copy()
$wsDates.Cells.Item(31, 4).PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --------------------------------------------------------------------------
# 3. "source" sheet - add the newly-cited source for the 2013 graphics
# --------------------------------------------------------------------------

$wsSource.Cells.Item(2, 1).Value = "Grafiken zur Studie 2013"
$wsSource.Cells.Item(2, 2).Value = "https://www.rosenheim24.de/deutschland/aengste-deutschen-2013-grafiken-studie-infocenters-rvversicherung-3094771.html"

# --------------------------------------------------------------------------
# 4. Window / selection state - "data" becomes the active tab
# --------------------------------------------------------------------------

$wsSource.Range("D5").Select()
$wsDates.Activate()
$wsDates.Range("H29").Select()

$wsData.Activate()
$wsData.Range("AA7").Select()

$wsDates.Range("E28").Select()
$wsSource.Range("B2").Select()

$wsData.Activate()
